$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.14295
$ws.Range("H2").Value = 0.2859
$ws.Range("I2").Value = 0.2060161421345645
$ws.Range("J2").Value = 0.15525926505184
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.0108025
$ws.Range("N2").Value = 0.021605
$ws.Range("Q2").Value = 0.001544217375
$ws.Range("R2").Value = 0.006176869499999999
$ws.Range("S2").Value = 0.2060161421345645
$ws.Range("T2").Value = 0.15525926505184
$ws.Range("I3").Value = 0.356759332696974
$ws.Range("J3").Value = 0.4032950371340628
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.0108025
$ws.Range("N3").Value = 0.021605
$ws.Range("Q3").Value = 0.002674130068333333
$ws.Range("R3").Value = 0.01604478041
$ws.Range("S3").Value = 0.356759332696974
$ws.Range("T3").Value = 0.4032950371340628
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07505833333333334
$ws.Range("H4").Value = 0.225175
$ws.Range("I4").Value = 0.1081722858928543
$ws.Range("J4").Value = 0.12228228404354
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.0108025
$ws.Range("N4").Value = 0.021605
$ws.Range("Q4").Value = 0.0008108176458333333
$ws.Range("R4").Value = 0.004864905875000001
$ws.Range("S4").Value = 0.1081722858928543
$ws.Range("T4").Value = 0.12228228404354
$ws.Range("G5").Value = 0.097247
$ws.Range("H5").Value = 0.194494
$ws.Range("I5").Value = 0.1401500648769499
$ws.Range("J5").Value = 0.1056208306995193
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.0108025
$ws.Range("N5").Value = 0.021605
$ws.Range("Q5").Value = 0.0010505107175
$ws.Range("R5").Value = 0.00420204287
$ws.Range("S5").Value = 0.1401500648769499
$ws.Range("T5").Value = 0.1056208306995193
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.02916366666666667
$ws.Range("H6").Value = 0.087491
$ws.Range("I6").Value = 0.04202998319108123
$ws.Range("J6").Value = 0.04751237621073989
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.0108025
$ws.Range("N6").Value = 0.021605
$ws.Range("Q6").Value = 0.0003150405091666667
$ws.Range("R6").Value = 0.001890243055
$ws.Range("S6").Value = 0.04202998319108123
$ws.Range("T6").Value = 0.04751237621073989
$ws.Range("G7").Value = 0.1019113333333333
$ws.Range("H7").Value = 0.305734
$ws.Range("I7").Value = 0.146872191207576
$ws.Range("J7").Value = 0.1660302068602982
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.0108025
$ws.Range("N7").Value = 0.021605
$ws.Range("Q7").Value = 0.001100897178333333
$ws.Range("R7").Value = 0.00660538307
$ws.Range("S7").Value = 0.146872191207576
$ws.Range("T7").Value = 0.1660302068602982
